$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated absenteeism records (rows 2-11), refactored for the ETL consolidator.
$data = @(
    @{Row=2;  A=89772; B="Emanuel Barbosa";          C="P&D";                     D="Doença";             E=8; F=45092; G=4209.49}
    @{Row=3;  A=14826; B="Kaique Alves";              C="Atendimento ao Cliente";  D="Consulta médica";    E=8; F=45094; G=2599.61}
    @{Row=4;  A=16361; B="Mariane Pinto";             C="Atendimento ao Cliente";  D="Viagem de negócios"; E=6; F=45094; G=5889.51}
    @{Row=5;  A=94445; B="Alexandre Caldeira";        C="P&D";                     D="Outros";             E=7; F=45090; G=9142.77}
    @{Row=6;  A=96531; B="João Guilherme da Costa";   C="Recursos Humanos";        D="Problemas pessoais"; E=3; F=45097; G=9054}
    @{Row=7;  A=23754; B="Lucas Gabriel Castro";      C="Recursos Humanos";        D="Viagem de negócios"; E=6; F=45086; G=5924.8}
    @{Row=8;  A=7956;  B="Dra. Stella Caldeira";      C="Operações";               D="Consulta médica";    E=5; F=45105; G=9516.66}
    @{Row=9;  A=41445; B="Letícia da Mota";           C="P&D";                     D="Doença";             E=2; F=45105; G=3755.18}
    @{Row=10; A=42193; B="Luiz Felipe Rezende";       C="Operações";               D="Consulta médica";    E=5; F=45083; G=12018.69}
    @{Row=11; A=19685; B="Vitor Hugo Martins";        C="Engenharia";              D="Doença";             E=1; F=45092; G=2720.22}
)

foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Cells.Item($r, 1).Value = $rec.A
    $ws.Cells.Item($r, 2).Value = $rec.B
    $ws.Cells.Item($r, 3).Value = $rec.C
    $ws.Cells.Item($r, 4).Value = $rec.D
    $ws.Cells.Item($r, 5).Value = $rec.E
    $ws.Cells.Item($r, 6).Value = $rec.F
    $ws.Cells.Item($r, 7).Value = $rec.G
}
